# Update "paises" COVID data sheet:
#  - refresh the "last updated" timestamp
#  - refresh per-country stats (totals/new/active/recovered/critical/deaths)
#  - a few countries changed rank order, so the country label in column A
#    for those rows now points at a different country than before
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 2 de Octubre de 2020 a las 23:17'
$ws.Range("B4").Value = 7535495
$ws.Range("C4").Value = 39975
$ws.Range("D4").Value = 4763732
$ws.Range("E4").Value = 2558418
$ws.Range("G4").Value = 685
$ws.Range("H4").Value = 213345
$ws.Range("A51").Value = 'Costa Rica'
$ws.Range("B51").Value = 77829
$ws.Range("C51").Value = 1001
$ws.Range("D51").Value = 42621
$ws.Range("E51").Value = 34278
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 930
$ws.Range("A52").Value = 'Honduras'
$ws.Range("B52").Value = 77598
$ws.Range("C52").Value = 698
$ws.Range("D52").Value = 28517
$ws.Range("E52").Value = 46701
$ws.Range("G52").Value = 27
$ws.Range("H52").Value = 2380
$ws.Range("A53").Value = 'Portugal'
$ws.Range("B53").Value = 77284
$ws.Range("C53").Value = 888
$ws.Range("D53").Value = 49359
$ws.Range("E53").Value = 25942
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 1983
$ws.Range("A54").Value = 'Etiopia'
$ws.Range("B54").Value = 76988
$ws.Range("C54").Value = 890
$ws.Range("D54").Value = 31677
$ws.Range("E54").Value = 44103
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 1208
$ws.Range("B66").Value = 46694
$ws.Range("C66").Value = 38
$ws.Range("D66").Value = 45945
$ws.Range("E66").Value = 448
$ws.Range("B85").Value = 19793
$ws.Range("C85").Value = 38
$ws.Range("D85").Value = 19357
$ws.Range("E85").Value = 316
$ws.Range("B101").Value = 11480
$ws.Range("C101").Value = 107
$ws.Range("D101").Value = 9351
$ws.Range("E101").Value = 2006
$ws.Range("B104").Value = 10398
$ws.Range("C104").Value = 44
$ws.Range("D104").Value = 9213
$ws.Range("E104").Value = 1151
$ws.Range("B113").Value = 7858
$ws.Range("C113").Value = 8
$ws.Range("D113").Value = 6322
$ws.Range("E113").Value = 1308
$ws.Range("B119").Value = 5783
$ws.Range("C119").Value = 4
$ws.Range("D119").Value = 4535
$ws.Range("E119").Value = 1069
$ws.Range("B122").Value = 5521
$ws.Range("C122").Value = 21
$ws.Range("D122").Value = 5061
$ws.Range("E122").Value = 349
$ws.Range("A124").Value = 'Angola'
$ws.Range("B124").Value = 5211
$ws.Range("C124").Value = 97
$ws.Range("D124").Value = 2215
$ws.Range("E124").Value = 2807
$ws.Range("G124").Value = 4
$ws.Range("H124").Value = 189
$ws.Range("A125").Value = 'Nicaragua'
$ws.Range("B125").Value = 5170
$ws.Range("D125").Value = 2913
$ws.Range("E125").Value = 2106
$ws.Range("H125").Value = 151
$ws.Range("B131").Value = 4847
$ws.Range("C131").Value = 4
$ws.Range("D131").Value = 3197
$ws.Range("E131").Value = 1621
$ws.Range("B134").Value = 4289
$ws.Range("C134").Value = 42
$ws.Range("D134").Value = 1130
$ws.Range("E134").Value = 2956
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 203
$ws.Range("B161").Value = 1818
$ws.Range("C161").Value = 9
$ws.Range("D161").Value = 1365
$ws.Range("E161").Value = 405
$ws.Range("B165").Value = 1346
$ws.Range("C165").Value = 3
$ws.Range("E165").Value = 43
$ws.Range("A181").Value = 'Curazao'
$ws.Range("B181").Value = 411
$ws.Range("C181").Value = 12
$ws.Range("D181").Value = 197
$ws.Range("E181").Value = 213
$ws.Range("H181").Value = 1
$ws.Range("A182").Value = 'San Martin (Parte Francesa)'
$ws.Range("B182").Value = 403
$ws.Range("D182").Value = 309
$ws.Range("E182").Value = 86
$ws.Range("H182").Value = 8
$ws.Range("B194").Value = 145
$ws.Range("C194").Value = 1
$ws.Range("E194").Value = 2
